$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the "Texas Notes" sheet entirely
$wb.Worksheets.Item("Texas Notes").Delete()

# 2. "About" sheet: point the hyperlink/text at the updated CEPE working-paper URL
$about = $wb.Worksheets.Item("About")
$newUrl = "https://ethz.ch/content/dam/ethz/special-interest/mtec/cepe/cepe-dam/documents/research/cepe-wp/CEPE_WP86.pdf"
$about.Range("B6").Value2 = $newUrl
$about.Hyperlinks.Delete()
$about.Hyperlinks.Add($about.Range("B6"), $newUrl)
$about.Range("B6").Style = "Hyperlink"

# 3. "MSCdtRPbQL" sheet: update the rebate-qualifying share number to 7.4%
$ml = $wb.Worksheets.Item("MSCdtRPbQL")
$ml.Range("C2").Value2 = 0.074
$ml.Range("A2").Select()

# 4. Make "About" the active tab/selection
$about.Activate()
$about.Range("A1").Select()
